# Applies: add a "2022-Q1" fund-holdings sheet (taking over the sheet
# slot previously used by "总计"), and re-create the "总计" (Total)
# summary sheet with the new 2022-Q1 row prepended to the old data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: the current last sheet ("总计") becomes the new "2022-Q1"
# sheet, keeping its original sheetId / relationship slot (matches the
# diff, where sheetId=6/rId6 goes from name="总计" to name="2022-Q1").
# ---------------------------------------------------------------------
$q1_2022 = $wb.Worksheets.Item("总计")
$q1_2022.Cells.Clear()
$q1_2022.Name = "2022-Q1"

# Use an existing, already-styled quarterly sheet as a formatting
# template so the re-used style index (s="2", bold/centered/bordered)
# stays identical instead of a brand new style being synthesized.
$template = $wb.Worksheets.Item("2021-Q4")

# Fund codes / scale / position figures are textual numeric-looking
# values in the source data (e.g. "000031" must keep its leading
# zero) -- force the target columns to text BEFORE any styled value
# is written so the numbers are not silently reinterpreted.
$q1_2022.Range("B2:G21").NumberFormat = "@"

# Copy the header row + the column-A "index" styling from the template
# sheet (this is the only styling the target layout actually uses).
$template.Range("B1:H1").Copy()
$q1_2022.Range("B1:H1").PasteSpecial(-4122)
$template.Range("A2").Copy()
$q1_2022.Range("A2:A21").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Header row
$q1_2022.Cells.Item(1,2).Value = "基金代码"
$q1_2022.Cells.Item(1,3).Value = "基金名称"
$q1_2022.Cells.Item(1,4).Value = "基金规模"
$q1_2022.Cells.Item(1,5).Value = "股票总仓位"
$q1_2022.Cells.Item(1,6).Value = "仓位占比"
$q1_2022.Cells.Item(1,7).Value = "持有市值(亿元)"
$q1_2022.Cells.Item(1,8).Value = "仓位排名"

$fundRows = @(
  @("501207", "华夏创新未来18个月封闭运作混合A", "67.75", "84.92", "6.96", "4.7154", 4),
  @("000031", "华夏复兴混合",                      "27.37", "89.15", "7.07", "1.9351", 4),
  @("007349", "华夏科技创新混合A",                 "14.63", "89.16", "7.04", "1.0300", 4),
  @("010518", "华夏先锋科技一年定期开放混合A",      "8.72",  "88.08", "7.46", "0.6505", 3),
  @("010106", "华夏核心科技6个月定期开放混合A",     "8.53",  "79.73", "6.79", "0.5792", 3),
  @("011184", "东方阿尔法招阳混合A",                "6.40",  "92.60", "8.88", "0.5683", 2),
  @("005358", "东方阿尔法精选灵活配置混合A",        "4.31",  "92.94", "7.87", "0.3392", 5),
  @("012568", "天弘高端制造混合型证券投资基金A",    "5.95",  "91.15", "4.21", "0.2505", 9),
  @("010519", "华夏先锋科技一年定期开放混合C",      "2.12",  "88.08", "7.46", "0.1582", 3),
  @("002746", "汇添富多策略定期开放灵活配置混合",   "4.55",  "64.75", "3.14", "0.1429", 8),
  @("010107", "华夏核心科技6个月定期开放混合C",     "1.32",  "79.73", "6.79", "0.0896", 3),
  @("007350", "华夏科技创新混合C",                  "1.08",  "89.16", "7.04", "0.0760", 4),
  @("005359", "东方阿尔法精选灵活配置混合C",        "0.54",  "92.94", "7.87", "0.0425", 5),
  @("012569", "天弘高端制造混合型证券投资基金C",    "0.82",  "91.15", "4.21", "0.0345", 9),
  @("011685", "创金合信先进装备股票A",              "0.73",  "92.01", "4.71", "0.0344", 8),
  @("004926", "中航军民融合精选混合A",              "0.35",  "91.27", "4.94", "0.0173", 7),
  @("004927", "中航军民融合精选混合C",              "0.27",  "91.27", "4.94", "0.0133", 7),
  @("003659", "山西证券策略精选灵活配置混合",       "0.31",  "84.52", "2.75", "0.0085", 10),
  @("011686", "创金合信先进装备股票C",              "0.17",  "92.01", "4.71", "0.0080", 8),
  @("011185", "东方阿尔法招阳混合C",                "0.08",  "92.60", "8.88", "0.0071", 2)
)

$r = 2
foreach ($row in $fundRows) {
  $q1_2022.Cells.Item($r, 1).Value = $r - 2
  $q1_2022.Cells.Item($r, 2).Value = $row[0]
  $q1_2022.Cells.Item($r, 3).Value = $row[1]
  $q1_2022.Cells.Item($r, 4).Value = $row[2]
  $q1_2022.Cells.Item($r, 5).Value = $row[3]
  $q1_2022.Cells.Item($r, 6).Value = $row[4]
  $q1_2022.Cells.Item($r, 7).Value = $row[5]
  $q1_2022.Cells.Item($r, 8).Value = $row[6]
  $r = $r + 1
}

# ---------------------------------------------------------------------
# Step 2: add a brand-new "总计" sheet right after "2022-Q1", re-stating
# the previous total/summary table plus the new 2022-Q1 row on top.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q1_2022)
$total.Name = "总计"

# Re-use the same template's header/index-column styling for the
# summary table -- this used to be exactly the "总计" layout.
$template.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$template.Range("A2").Copy()
$total.Range("A2:A7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$total.Cells.Item(1,2).Value = "日期"
$total.Cells.Item(1,3).Value = "持有数量(只)"
$total.Cells.Item(1,4).Value = "持有市值(亿元)"

$summaryRows = @(
  @("2022-Q1", 20, 10.7),
  @("2021-Q4", 24, 29.51),
  @("2021-Q3", 28, 28.31),
  @("2021-Q2", 27, 27.12),
  @("2021-Q1", 22, 12.08),
  @("2020-Q4", 12, 15.1)
)

$r = 2
foreach ($row in $summaryRows) {
  $total.Cells.Item($r, 1).Value = $r - 2
  $total.Cells.Item($r, 2).Value = $row[0]
  $total.Cells.Item($r, 3).Value = $row[1]
  $total.Cells.Item($r, 4).Value = $row[2]
  $r = $r + 1
}

$q1_2022.Select()
